$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.729.08"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.100.10"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "'343.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "'0.5193"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("D8").Value = "'0.4387"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "'53.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "
$ws.Range("D10").Value = "'0.09239"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.88%  "
$ws.Range("D11").Value = "'1.168"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").Value = "'24.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.37%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "2.133.78"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.796"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "'8.174"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "'102.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("D17").Value = "'0.00001154"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "'21.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "'0.06667"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "'6.218"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").Value = "29.738.97"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("D24").Value = "'12.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").Value = "'2.306"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").Value = "2.308.19"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "'21.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'162.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'2.500"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.59%  "
$ws.Range("D30").Value = "'133.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").Value = "'1.131"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.08%  "
$ws.Range("D32").Value = "'1.701"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("D33").Value = "'0.1052"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "'6.201"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "'3.954"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "'6.359"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.03%  "
$ws.Range("D37").Value = "'10.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "'0.06724"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("D40").Value = "'0.6987"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.35%  "
$ws.Range("D41").Value = "'12.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").Value = "'1.329"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("D44").Value = "'0.6797"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.51%  "
$ws.Range("D45").Value = "'14.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").Value = "'2.334"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.04%  "
$ws.Range("D47").Value = "'0.00000000358"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").Value = "'3.621"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").Value = "'1.219"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'1.201"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("D51").Value = "'81.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.35%  "
